$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.326392292976379
$ws.Range("B1").Value = 3.291231870651245
$ws.Range("C1").Value = 5.592843532562256
$ws.Range("D1").Value = 1.709142446517944
$ws.Range("E1").Value = 0.9997802376747131
